$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column B for rows 9-11
$ws.Range("B9").Value = -5
$ws.Range("B10").Value = -5
$ws.Range("B11").Value = -5

# Move the active selection/cell to F10
$ws.Range("F10").Select()
